$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Normalize study-name text that previously used Unicode hyphens (U+2010) and a
# trailing non-breaking space into plain ASCII hyphens (part of extending the
# resistance evidence base to grazoprevir / C-EDGE & C-SCAPE trials).
$ws.Range("A66").Value = "C-SCAPE"
$ws.Range("B66").Value = "C-SCAPE"

$ws.Range("A67").Value = "C-EDGE IBLD"
$ws.Range("B67").Value = "C-EDGE IBLD"

$ws.Range("A68").Value = "C-EDGE Head-2-head"
$ws.Range("B68").Value = "C-EDGE Head-2-head"

# Select the whole sheet (matches the saved view state in the workbook).
$ws.Cells.Select()
